# Added electricity consumption profile generator.
# Variable appliances still need to be added.

$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet from "Sheet1" to "Materials" ---
$materials = $wb.Worksheets.Item("Sheet1")
$materials.Name = "Materials"

# --- Add a new worksheet named "Inputs" right after "Materials" ---
$inputs = $wb.Worksheets.Add($null, $materials)
$inputs.Name = "Inputs"

# --- Materials sheet view tweaks ---
$materials.Select()
$materials.Columns.Item(1).AutoFit()
$materials.Range("A7").Select()

# Zoom the Materials sheet window to 220%
$excel.ActiveWindow.Zoom = 220

# --- Make Inputs the active (selected) sheet/tab ---
$inputs.Select()

$excel.ActiveWindow.WindowState = -4140  # xlMinimized
